$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a range to Text type without leaving a residual style
# change (NumberFormat "@" forces literal-text parsing of numeric-looking
# strings; resetting Style back to "Normal" drops the cell back onto the
# original, unstyled xf so only the VALUE differs from the source file).
function Set-TextValue($rangeAddr, $text) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '42.891.39'
$ws.Range('E2').Value = '  -0.08%  '

$ws.Range('D3').Value = '2.312.23'
$ws.Range('E3').Value = '  +0.32%  '

Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.03%  '

Set-TextValue 'D5' '301.62'
$ws.Range('E5').Value = '  -1.79%  '

Set-TextValue 'D6' '96.62'
$ws.Range('E6').Value = '  -0.89%  '

Set-TextValue 'D7' '0.504'
$ws.Range('E7').Value = '  -1.40%  '

$ws.Range('E8').Value = '  +0.13%  '

Set-TextValue 'D9' '0.497'
$ws.Range('E9').Value = '  -1.69%  '

Set-TextValue 'D10' '34.94'
$ws.Range('E10').Value = '  -2.23%  '

Set-TextValue 'D11' '19.32'
$ws.Range('E11').Value = '  +6.08%  '

Set-TextValue 'D12' '0.0789'
$ws.Range('E12').Value = '  -0.08%  '

$ws.Range('E13').Value = '  +0.22%  '

Set-TextValue 'D14' '6.84'
$ws.Range('E14').Value = '  +0.83%  '

$ws.Range('D15').Value = '2.674.81'
$ws.Range('E15').Value = '  +0.60%  '

$ws.Range('D16').Value = '2.308.51'
$ws.Range('E16').Value = '  +0.30%  '

Set-TextValue 'D17' '0.784'
$ws.Range('E17').Value = '  -0.27%  '

$ws.Range('D18').Value = '42.871.40'
$ws.Range('E18').Value = '  +0.04%  '

Set-TextValue 'D19' '12.43'
$ws.Range('E19').Value = '  -2.24%  '

$ws.Range('D20').Value = '0.0₃0889'
$ws.Range('E20').Value = '  -1.53%  '

Set-TextValue 'D21' '6.02'
$ws.Range('E21').Value = '  -0.33%  '

Set-TextValue 'D22' '67.57'
$ws.Range('E22').Value = '  -0.38%  '

Set-TextValue 'D23' '235.62'
$ws.Range('E23').Value = '  -0.42%  '

Set-TextValue 'D24' '2.22'
$ws.Range('E24').Value = '  +3.32%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D25' '1.00'
$ws.Range('E25').Value = '  +0.05%  '

$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D26' '2.43'
$ws.Range('E26').Value = '  -2.03%  '

Set-TextValue 'D27' '24.67'
$ws.Range('E27').Value = '  -3.05%  '

$ws.Range('E28').Value = '  +0.33%  '

Set-TextValue 'D29' '163.83'
$ws.Range('E29').Value = '  -1.90%  '

Set-TextValue 'D30' '9.05'
$ws.Range('E30').Value = '  -0.08%  '

Set-TextValue 'D31' '32.52'
$ws.Range('E31').Value = '  -1.90%  '

Set-TextValue 'D32' '1.00'
$ws.Range('E32').Value = '  +0.06%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '4.97'
$ws.Range('E33').Value = '  -0.97%  '

$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D34' '17.73'
$ws.Range('E34').Value = '  +2.80%  '

Set-TextValue 'D35' '4.46'
$ws.Range('E35').Value = '  -7.22%  '

$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D36' '2.35'
$ws.Range('E36').Value = '  -1.17%  '

$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D37' '0.0696'
$ws.Range('E37').Value = '  +0.63%  '

$ws.Range('E38').Value = '  -0.98%  '

Set-TextValue 'D39' '1.76'
$ws.Range('E39').Value = '  -0.13%  '

Set-TextValue 'D40' '2.75'
$ws.Range('E40').Value = '  +0.54%  '

Set-TextValue 'D41' '0.108'
$ws.Range('E41').Value = '  -1.19%  '

$ws.Range('D42').Value = '1.975.37'
$ws.Range('E42').Value = '  -1.54%  '

Set-TextValue 'D43' '10.58'
$ws.Range('E43').Value = '  +5.58%  '

Set-TextValue 'D44' '18.61'
$ws.Range('E44').Value = '  +3.58%  '

Set-TextValue 'D45' '0.0279'
$ws.Range('E45').Value = '  -1.22%  '

$ws.Range('E46').Value = '  -4.35%  '

Set-TextValue 'D47' '2.76'
$ws.Range('E47').Value = '  -1.08%  '

$ws.Range('D48').Value = '2.539.31'
$ws.Range('E48').Value = '  +0.49%  '

$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D49' '2.85'
$ws.Range('E49').Value = '  +0.13%  '

$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D50' '53.27'
$ws.Range('E50').Value = '  -1.16%  '

Set-TextValue 'D51' '71.92'
$ws.Range('E51').Value = '  -0.12%  '
